# Update 2p0. Convention change to support multi-axle vehicles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet renamed to drop the redundant "Driveline" prefix now that the
# folder/axle naming convention changed (Axle_1 -> Axle1).
$ws.Name = "Axle1_None"

# Narrow column B very slightly to match the new layout.
$ws.Columns.Item(2).ColumnWidth = 11.8333333

# Move the active selection on the frozen (bottom-right) pane.
$ws.Range("H5").Select()
